# Excel import -> improved email validation
# Populate the "users" import-template example sheet (rows 2-17) with
# sample data exercising multi-column sort / duplicate detection, and
# clear out the old placeholder rows 9-17 (which only had a formatted,
# empty D cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sample data rows (A:E), for rows 2..17 -------------------------------
$rows = @(
  @("a.multi.column.sort@example.com","a_test","AAAA","AAAB","AAAA"),
  @("b.multi.column.sort@example.com","b_test","AAAB","AAAA","AAAA"),
  @("c.multi.column.sort@example.com","c_test","AAAA","AAAC","AAAD"),
  @("d.multi.column.sort@example.com","d_test","AAAC","AAAA","AAAD"),
  @("e.multi.column.sort@example.com","e_test","AAAA","AAAB","AAAC"),
  @("f.multi.column.sort@example.com","f_test","AAAB","AAAA","AAAC"),
  @("g.multi.column.sort@example.com","g_test","AAAA","AAAC","AAAB"),
  @("h.multi.column.sort@example.com","h_test","AAAC","AAAA","AAAB"),
  @("i.multi.column.sort@example.com","I_test","AAAA","AAAB","AAAA"),
  @("j.multi.column.sort@example.com","j_test","AAAB","AAAA","AAAA"),
  @("k.multi.column.sort@example.com","k_test","AAAA","AAAC","AAAD"),
  @("l.multi.column.sort@example.com","l_test","AAAC","AAAA","AAAD"),
  @("m.multi.column.sort@example.com","m_test","AAAA","AAAB","AAAC"),
  @("n.multi.column.sort@example.com","n_test","AAAB","AAAA","AAAC"),
  @("o.multi.column.sort@example.com","o_test","AAAA","AAAC","AAAB"),
  @("p.multi.column.sort@example.com","p_test","AAAC","AAAA","AAAB")
)

$rowCount = $rows.Count
$colCount = 5

$data = New-Object 'object[,]' $rowCount,$colCount
for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $data[$r,$c] = $rows[$r][$c]
    }
}

$lastRow = 1 + $rowCount
$target = $ws.Range("A2:E$lastRow")
$target.Value2 = $data

# give the new data rows the same row height Excel/LO used for rows 3-17
for ($r = 3; $r -le $lastRow; $r++) {
    $ws.Rows.Item($r).RowHeight = 16.4
}

# widen column A slightly (closest achievable match to the target 31.49 chars,
# since this runtime snaps widths to 1/6-character increments)
$ws.Columns.Item(1).ColumnWidth = 30.666666666666668

# update the saved selection / active cell shown when the sheet is reopened
$ws.Range("E10:E17").Select() | Out-Null
